$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 777.4
$ws.Range("I29").Value = 596.75
$ws.Range("J29").Value = 1500
$ws.Range("K29").Value = 1790.25
$ws.Range("L29").Value = 4500
$ws.Range("M29").Value = -1509.25
$ws.Range("N29").Value = -5062
$ws.Range("H38").Value = 996.6429000000001
$ws.Range("I38").Value = 65.3
$ws.Range("K38").Value = 195.9
$ws.Range("M38").Value = 176.1
$ws.Range("H43").Value = 536.5833
$ws.Range("I43").Value = 628.5714
$ws.Range("J43").Value = 407.8
$ws.Range("K43").Value = 628.5714
$ws.Range("L43").Value = 407.8
$ws.Range("M43").Value = -559.5714
$ws.Range("N43").Value = -545.8
$ws.Range("H58").Value = 1529.7727
$ws.Range("I58").Value = 114
$ws.Range("J58").Value = 2945.5454
$ws.Range("K58").Value = 342
$ws.Range("L58").Value = 8836.636200000001
$ws.Range("M58").Value = -192
$ws.Range("N58").Value = -9136.636200000001
$ws.Range("H129").Value = 1109.7
$ws.Range("I129").Value = 649.25
$ws.Range("J129").Value = 1224.8125
$ws.Range("K129").Value = 1947.75
$ws.Range("L129").Value = 3674.4375
$ws.Range("M129").Value = 3052.25
$ws.Range("N129").Value = -13674.4375
$ws.Range("H137").Value = 1489.65
$ws.Range("I137").Value = 1435.6086
$ws.Range("J137").Value = 1562.7646
$ws.Range("K137").Value = 4306.825800000001
$ws.Range("L137").Value = 4688.293799999999
$ws.Range("M137").Value = -1756.825800000001
$ws.Range("N137").Value = -9788.293799999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1422.7354
$ws.Range("I2").Value = 1041.5454
$ws.Range("J2").Value = 2121.5833
$ws.Range("K2").Value = 1041.5454
$ws.Range("L2").Value = 2121.5833
$ws.Range("M2").Value = -928.5454
$ws.Range("N2").Value = -2347.5833
$ws.Range("H94").Value = 28000
$ws.Range("J94").Value = 28000
$ws.Range("L94").Value = 28000
$ws.Range("N94").Value = -29802
$ws.Range("H116").Value = 1422.7354
$ws.Range("I116").Value = 1041.5454
$ws.Range("J116").Value = 2121.5833
$ws.Range("K116").Value = 1041.5454
$ws.Range("L116").Value = 2121.5833
$ws.Range("M116").Value = 1252.4546
$ws.Range("N116").Value = -6709.5833

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1422.7354
$ws.Range("I3").Value = 1041.5454
$ws.Range("J3").Value = 2121.5833
$ws.Range("K3").Value = 1041.5454
$ws.Range("L3").Value = 2121.5833
$ws.Range("M3").Value = -927.5454
$ws.Range("N3").Value = -2349.5833
$ws.Range("H20").Value = 43026.12
$ws.Range("I20").Value = 54655.895
$ws.Range("K20").Value = 54655.895
$ws.Range("M20").Value = -54408.895

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 17001.666
$ws.Range("J3").Value = 3
$ws.Range("L3").Value = 3
$ws.Range("N3").Value = -229
$ws.Range("H16").Value = 1563.25
$ws.Range("I16").Value = 1417.6666
$ws.Range("J16").Value = 2000
$ws.Range("K16").Value = 1417.6666
$ws.Range("L16").Value = 2000
$ws.Range("M16").Value = -1130.6666
$ws.Range("N16").Value = -2574
$ws.Range("H113").Value = 1563.25
$ws.Range("I113").Value = 1417.6666
$ws.Range("J113").Value = 2000
$ws.Range("K113").Value = 1417.6666
$ws.Range("L113").Value = 2000
$ws.Range("M113").Value = 752.3334
$ws.Range("N113").Value = -6340
$ws.Range("H132").Value = 4333.1113
$ws.Range("I132").Value = 4000
$ws.Range("J132").Value = 4749.5
$ws.Range("K132").Value = 12000
$ws.Range("L132").Value = 14248.5
$ws.Range("M132").Value = -9470
$ws.Range("N132").Value = -19308.5
$ws.Range("H141").Value = 47479.5
$ws.Range("I141").Value = 40000
$ws.Range("J141").Value = 49972.668
$ws.Range("K141").Value = 40000
$ws.Range("L141").Value = 49972.668
$ws.Range("M141").Value = -34820
$ws.Range("N141").Value = -60332.668

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 765.8148
$ws.Range("I122").Value = 498.66666
$ws.Range("J122").Value = 1300.1111
$ws.Range("K122").Value = 4487.99994
$ws.Range("L122").Value = 11700.9999
$ws.Range("M122").Value = -2037.99994
$ws.Range("N122").Value = -16600.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3462.647
$ws.Range("I132").Value = 2665.1667
$ws.Range("J132").Value = 3897.6365
$ws.Range("K132").Value = 7995.500100000001
$ws.Range("L132").Value = 11692.9095
$ws.Range("M132").Value = -5465.500100000001
$ws.Range("N132").Value = -16752.9095

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H5").Value = 15787.143
$ws.Range("I5").Value = 18000
$ws.Range("J5").Value = 14127.5
$ws.Range("K5").Value = 18000
$ws.Range("L5").Value = 14127.5
$ws.Range("M5").Value = -17887
$ws.Range("N5").Value = -14353.5
$ws.Range("H61").Value = 16792.857
$ws.Range("I61").Value = 20372.818
$ws.Range("J61").Value = 3666.3333
$ws.Range("K61").Value = 20372.818
$ws.Range("L61").Value = 3666.3333
$ws.Range("M61").Value = -20170.818
$ws.Range("N61").Value = -4070.3333
$ws.Range("H113").Value = 16792.857
$ws.Range("I113").Value = 20372.818
$ws.Range("J113").Value = 3666.3333
$ws.Range("K113").Value = 20372.818
$ws.Range("L113").Value = 3666.3333
$ws.Range("M113").Value = -18202.818
$ws.Range("N113").Value = -8006.3333
$ws.Range("H132").Value = 6690.174
$ws.Range("I132").Value = 6829.1665
$ws.Range("K132").Value = 20487.4995
$ws.Range("M132").Value = -17957.4995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 21474.5
$ws.Range("J21").Value = 11966
$ws.Range("L21").Value = 11966
$ws.Range("N21").Value = -12436
$ws.Range("H25").Value = 29970.8
$ws.Range("J25").Value = 29970.8
$ws.Range("L25").Value = 29970.8
$ws.Range("N25").Value = -30556.8
$ws.Range("H33").Value = 15000
$ws.Range("J33").Value = 15000
$ws.Range("L33").Value = 15000
$ws.Range("N33").Value = -15500
$ws.Range("H35").Value = 21474.5
$ws.Range("J35").Value = 11966
$ws.Range("L35").Value = 11966
$ws.Range("N35").Value = -12546
$ws.Range("H36").Value = 15000
$ws.Range("J36").Value = 15000
$ws.Range("L36").Value = 15000
$ws.Range("N36").Value = -15500
$ws.Range("H37").Value = 50009.668
$ws.Range("J37").Value = 50014.5
$ws.Range("L37").Value = 50014.5
$ws.Range("N37").Value = -50420.5
$ws.Range("H122").Value = 10138.75
$ws.Range("I122").Value = 10000
$ws.Range("K122").Value = 30000
$ws.Range("M122").Value = -27550
